# Auto-generated edit script applying numeric updates to the Brynhildr_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 59536.23
$ws.Range("I11").Value = 59536.23
$ws.Range("K11").Value = 59536.23
$ws.Range("M11").Value = -59396.23

$ws.Range("H39").Value = 1075.3334
$ws.Range("J39").Value = 1162.6666
$ws.Range("L39").Value = 3487.9998
$ws.Range("N39").Value = -4079.9998

$ws.Range("H86").Value = 7581.6
$ws.Range("I86").Value = 7938.6
$ws.Range("K86").Value = 7938.6
$ws.Range("M86").Value = -6815.6

$ws.Range("H89").Value = 7581.6
$ws.Range("I89").Value = 7938.6
$ws.Range("K89").Value = 39693
$ws.Range("M89").Value = -34077

$ws.Range("H101").Value = 278
$ws.Range("I101").Value = 278
$ws.Range("K101").Value = 834
$ws.Range("M101").Value = 788

$ws.Range("H103").Value = 1820
$ws.Range("I103").Value = 1500
$ws.Range("J103").Value = 2033.3334
$ws.Range("K103").Value = 4500
$ws.Range("L103").Value = 6100.0002
$ws.Range("M103").Value = -3914
$ws.Range("N103").Value = -7272.0002

$ws.Range("H135").Value = 3395
$ws.Range("I135").Value = 741.2
$ws.Range("K135").Value = 6670.8
$ws.Range("M135").Value = -4135.8

$ws.Range("H137").Value = 27785768
$ws.Range("I137").Value = 38463796
$ws.Range("J137").Value = 22900
$ws.Range("K137").Value = 115391388
$ws.Range("L137").Value = 68700
$ws.Range("M137").Value = -115388838
$ws.Range("N137").Value = -73800

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1962.1
$ws.Range("I45").Value = 1152.75
$ws.Range("J45").Value = 5199.5
$ws.Range("K45").Value = 1152.75
$ws.Range("L45").Value = 5199.5
$ws.Range("M45").Value = -775.75
$ws.Range("N45").Value = -5953.5

$ws.Range("H61").Value = 5558023.5
$ws.Range("I61").Value = 2613.1765
$ws.Range("K61").Value = 2613.1765
$ws.Range("M61").Value = -2401.1765

$ws.Range("H74").Value = 782146.5
$ws.Range("I74").Value = 872788.4
$ws.Range("K74").Value = 872788.4
$ws.Range("M74").Value = -871914.4

$ws.Range("H77").Value = 782146.5
$ws.Range("I77").Value = 872788.4
$ws.Range("K77").Value = 4363942
$ws.Range("M77").Value = -4359574

$ws.Range("H110").Value = 1875.5
$ws.Range("I110").Value = 1499
$ws.Range("K110").Value = 1499
$ws.Range("M110").Value = 546

$ws.Range("H122").Value = 1835.2222
$ws.Range("I122").Value = 1707.25
$ws.Range("K122").Value = 5121.75
$ws.Range("M122").Value = -2671.75

$ws.Range("H132").Value = 4155.9414
$ws.Range("I132").Value = 2168.5454
$ws.Range("K132").Value = 6505.6362
$ws.Range("M132").Value = -3975.6362

$ws.Range("H136").Value = 5558023.5
$ws.Range("I136").Value = 2613.1765
$ws.Range("K136").Value = 7839.529500000001
$ws.Range("M136").Value = -5289.529500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3969700
$ws.Range("I134").Value = 1387.4736
$ws.Range("K134").Value = 4162.4208
$ws.Range("M134").Value = -1627.4208

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1593970.9
$ws.Range("I31").Value = 2460183.8
$ws.Range("J31").Value = 5914
$ws.Range("K31").Value = 2460183.8
$ws.Range("L31").Value = 5914
$ws.Range("M31").Value = -2459888.8
$ws.Range("N31").Value = -6504

$ws.Range("H34").Value = 1593970.9
$ws.Range("I34").Value = 2460183.8
$ws.Range("J34").Value = 5914
$ws.Range("K34").Value = 2460183.8
$ws.Range("L34").Value = 5914
$ws.Range("M34").Value = -2459981.8
$ws.Range("N34").Value = -6318

$ws.Range("H107").Value = 439.2
$ws.Range("J107").Value = 672.3333
$ws.Range("L107").Value = 672.3333
$ws.Range("N107").Value = -4512.3333

$ws.Range("H134").Value = 4311.517
$ws.Range("I134").Value = 2524.7407
$ws.Range("K134").Value = 7574.222099999999
$ws.Range("M134").Value = -5039.222099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 82.333336
$ws.Range("I2").Value = 27.11111
$ws.Range("J2").Value = 115.46667
$ws.Range("K2").Value = 162.66666
$ws.Range("L2").Value = 692.8000199999999
$ws.Range("M2").Value = -49.66666000000001
$ws.Range("N2").Value = -918.8000199999999

$ws.Range("H5").Value = 5361709
$ws.Range("I5").Value = 3968887.5
$ws.Range("K5").Value = 11906662.5
$ws.Range("M5").Value = -11906550.5

$ws.Range("H82").Value = 12500
$ws.Range("I82").Value = 7500
$ws.Range("K82").Value = 22500
$ws.Range("M82").Value = -22094

$ws.Range("H85").Value = 12500
$ws.Range("I85").Value = 7500
$ws.Range("K85").Value = 22500
$ws.Range("M85").Value = -21096

$ws.Range("H131").Value = 6392.4
$ws.Range("J131").Value = 9315.315000000001
$ws.Range("L131").Value = 27945.945
$ws.Range("N131").Value = -38025.945

$ws.Range("H132").Value = 1619.6
$ws.Range("I132").Value = 1224.75
$ws.Range("K132").Value = 11022.75
$ws.Range("M132").Value = -8492.75

$ws.Range("H135").Value = 5361709
$ws.Range("I135").Value = 3968887.5
$ws.Range("K135").Value = 35719987.5
$ws.Range("M135").Value = -35717452.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1059.2
$ws.Range("I102").Value = 954.6667
$ws.Range("K102").Value = 954.6667
$ws.Range("M102").Value = 667.3333

$ws.Range("H113").Value = 1329.375
$ws.Range("I113").Value = 1329.375
$ws.Range("K113").Value = 1329.375
$ws.Range("M113").Value = 840.625

$ws.Range("H122").Value = 28995.406
$ws.Range("I122").Value = 31411.03
$ws.Range("K122").Value = 94233.09
$ws.Range("M122").Value = -91783.09

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2319.8
$ws.Range("I7").Value = 2319.8
$ws.Range("K7").Value = 2319.8
$ws.Range("M7").Value = -2207.8

$ws.Range("H93").Value = 5951.643
$ws.Range("I93").Value = 1811.875
$ws.Range("J93").Value = 11471.333
$ws.Range("K93").Value = 1811.875
$ws.Range("L93").Value = 11471.333
$ws.Range("M93").Value = -563.875
$ws.Range("N93").Value = -13967.333

$ws.Range("H126").Value = 2319.8
$ws.Range("I126").Value = 2319.8
$ws.Range("K126").Value = 6959.400000000001
$ws.Range("M126").Value = -4489.400000000001

$ws.Range("H132").Value = 1452423.9
$ws.Range("I132").Value = 3032963.5
$ws.Range("K132").Value = 9098890.5
$ws.Range("M132").Value = -9096360.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1383.25
$ws.Range("I96").Value = 1430.4286
$ws.Range("J96").Value = 1346.5555
$ws.Range("K96").Value = 1430.4286
$ws.Range("L96").Value = 1346.5555
$ws.Range("M96").Value = -57.42859999999996
$ws.Range("N96").Value = -4092.5555

$ws.Range("H126").Value = 1314
$ws.Range("I126").Value = 1087.6666
$ws.Range("J126").Value = 1766.6666
$ws.Range("K126").Value = 3262.9998
$ws.Range("L126").Value = 5299.9998
$ws.Range("M126").Value = -792.9998000000001
$ws.Range("N126").Value = -10239.9998

$ws.Range("H132").Value = 3403782.5
$ws.Range("I132").Value = 3548220
$ws.Range("K132").Value = 10644660
$ws.Range("M132").Value = -10642130
